$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column letter -> new value, applying the weekly re-sync of dates/volumes/prices.
$updates = @{
    2  = @{ D = 44418; M = 40 }
    3  = @{ D = 44432; M = 30 }
    4  = @{ D = 44435; M = 130; N = 1300; O = 1300; P = 1300; S = 1300 }
    5  = @{ D = 44424; M = 50;  N = 1200; O = 1200; P = 1200; S = 1200 }
    6  = @{ D = 44417; M = 80;  N = 1200; O = 1200; P = 1200; S = 1200 }
    7  = @{ D = 44438; M = 60 }
    8  = @{ D = 44431; M = 100; N = 1300; O = 1300; P = 1300; S = 1300 }
    9  = @{ D = 44357; M = 35;  N = 1000; O = 1000; P = 1000; S = 1000 }
    10 = @{ D = 44343; M = 60;  N = 1300; O = 1300; P = 1300; S = 1300 }
    11 = @{ D = 44405; M = 50;  N = 1200; O = 1200; P = 1200; S = 1200 }
    12 = @{ D = 44476; M = 80 }
    13 = @{ D = 44473; M = 120; N = 1200; O = 1200; P = 1200; S = 1200 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $cellRef = "$col$row"
        $ws.Range($cellRef).Value = $cols[$col]
    }
}
